$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (the source data keeps these as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '70.762.58'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '3.803.75'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '704.15'
$ws.Range("E5").Value = '  +2.21%  '
$ws.Range("D6").Value = '169.20'
$ws.Range("E6").Value = '  -2.24%  '
$ws.Range("D7").Value = '3.798.04'
$ws.Range("E7").Value = '  -1.33%  '
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").Value = '0.455'
$ws.Range("E12").Value = '  -1.44%  '
$ws.Range("D13").Value = '0.0000252'
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("D14").Value = '36.15'
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").Value = '4.452.41'
$ws.Range("E15").Value = '  -0.97%  '
$ws.Range("D16").Value = '3.815.07'
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = '70.892.28'
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '0.115'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '7.15'
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("D20").Value = '17.21'
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").Value = '491.24'
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").Value = '10.54'
$ws.Range("E22").Value = '  -4.87%  '
$ws.Range("D23").Value = '0.726'
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("D24").Value = '85.45'
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = '12.04'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '10.45'
$ws.Range("E27").Value = '  -0.84%  '
$ws.Range("D28").Value = '3.963.05'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("D31").Value = '3.07'
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("D32").Value = '7.34'
$ws.Range("E32").Value = '  -3.89%  '
$ws.Range("E33").Value = '  -4.09%  '
$ws.Range("D34").Value = '29.10'
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("E35").Value = '  -4.31%  '
$ws.Range("D36").Value = '3.778.62'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '9.03'
$ws.Range("E38").Value = '  -2.68%  '
$ws.Range("E39").Value = '  -2.65%  '
$ws.Range("E40").Value = '  +2.52%  '
$ws.Range("E41").Value = '  -3.56%  '
$ws.Range("D42").Value = '5.91'
$ws.Range("E42").Value = '  -2.54%  '
$ws.Range("D43").Value = '3.27'
$ws.Range("E43").Value = '  -5.01%  '
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("E45").Value = '  +0.16%  '
$ws.Range("D46").Value = '164.12'
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").Value = '0.000307'
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").Value = '48.68'
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").Value = '422.22'
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("D50").Value = '8.68'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").Value = '0.293'
$ws.Range("E51").Value = '  -3.16%  '
